$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

$ws.Range("G4").Value = "1"
$ws.Range("G5").Value = "21"

$ws.Range("G5").Select()
